# Apply the "Card20" row 33 (sheet row 35) reset edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Clear A35 (was "20")
$ws.Range("A35").Value = ""

# Fill B35:K35 with "nan"
$ws.Range("B35:K35").Value = "nan"

# Clear L35:N35 (Date, Event, Correction)
$ws.Range("L35:N35").Value = ""

# O35 (Serviced by) remains unchanged: "م/محمد**محمود ايهاب**ابراهيم"
$ws.Range("O35").Value = "م/محمد**محمود ايهاب**ابراهيم"
